# Add files via upload
# Sprint 1 Problem Definition: insert a new set of 4 engineering tasks
# ("No" task-complete markers included) for the "Presentation Outline"
# story, directly above its existing row, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1 Problem Definition")

# Make 9 new rows of room (rows 7-15) above the old row 7 ("Bruce" /
# "Presentation Outline" story), which pushes everything from the old
# row 7 onward down by 9 rows (8 new content rows + 1 blank spacer row,
# matching the existing blank-row-between-stories pattern used
# throughout the sheet).
$ws.Rows("7:15").Insert()

# New engineering task rows for the story that is now anchored at A16/B16.
$ws.Range("C7").Value = "1. Identify terms for analysis, e.g. bias, fact, opinion, pedigree of source"
$ws.Range("E8").Value = "No"

$ws.Range("C9").Value = "2. Establish weightings and scales for these terms, e.g. 1-10 bias, but 1-100 for facts"
$ws.Range("E10").Value = "No"

$ws.Range("C11").Value = "3. Identify if these terms are immutable or user preference, e.g. where on the bias scale am I, but facts are facts"
$ws.Range("E12").Value = "No"

$ws.Range("C13").Value = "4. Identify current research to support the above items. Look at fraud and deception as similar topics."
$ws.Range("E14").Value = "No"

# Row 15 stays blank (spacer row before the next story), matching the
# layout used elsewhere in the sheet.

# Move/restore the active selection to C14, matching the saved view state.
$ws.Range("C14").Select() | Out-Null
